# Updated symbol list on Sun Dec 18 05:26:43 UTC 2022 with GitHub Actions
#
# The "Price" column (D) stores values as text, and several of the new
# prices still look numeric, so a leading apostrophe is used to force
# Excel to keep them as text (matching the existing cell type) instead
# of silently re-typing the cell as a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates -------------------------------------------------
$ws.Range("D2").Value  = "'248.57"
$ws.Range("D3").Value  = "'22.46"
$ws.Range("D4").Value  = "'5.636"
$ws.Range("D5").Value  = "'0.05612"
$ws.Range("D6").Value  = "'3.397"
$ws.Range("D7").Value  = "'6.473"
$ws.Range("D8").Value  = "'1.078"
$ws.Range("D9").Value  = "'0.8039"
$ws.Range("D10").Value = "'0.1430"
$ws.Range("D11").Value = "'0.07437"
$ws.Range("D12").Value = "'0.03189"
$ws.Range("D13").Value = "'0.02992"
$ws.Range("D14").Value = "'0.09265"
$ws.Range("D15").Value = "'0.001665"
$ws.Range("D16").Value = "'3.249"
$ws.Range("D17").Value = "'0.04747"
$ws.Range("D18").Value = "'0.0005739"
$ws.Range("D19").Value = "'0.006261"
$ws.Range("D20").Value = "'0.001053"
$ws.Range("D21").Value = "'0.003817"
$ws.Range("D25").Value = "'2.118"
$ws.Range("D27").Value = "'0.1276"
$ws.Range("D40").Value = "'0.04192"

# --- Row 18: "One" (ONE) lost its "Worst in 24h" badge -------------------------
$ws.Range("E18").Value = "17OneONE"

# --- Rows 41-43: coin rotation (KickToken / BKEXToken / CEJI reshuffled) ------
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1048"
$ws.Range("E41").Value = "40BKEXTokenBKK"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.002969"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.003254"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"

# --- Remaining price updates ---------------------------------------------------
$ws.Range("D44").Value = "'0.009827"
$ws.Range("D45").Value = "'0.00005664"
$ws.Range("D48").Value = "'0.02944"
$ws.Range("D49").Value = "'0.00002099"
